$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 11 and 12 swap position: BinanceUSD now ranks above Dogecoin
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"

# Updated Price (D) and Volume(1h) (E) figures for every listed coin.
# A leading apostrophe forces Excel to keep number-looking text (e.g. "1.001")
# stored as text instead of silently converting it to a numeric value.
$priceVolume = @{
  2 = @("23.223.15", "  +0.87%  ")
  3 = @("1.605.20", "  +0.40%  ")
  4 = @("'1.001", "  -0.01%  ")
  5 = @("'1.001", "  -0.03%  ")
  6 = @("'304.17", "  +0.81%  ")
  7 = @("'0.3764", "  -0.40%  ")
  8 = @("'52.33", "  +5.28%  ")
  9 = @("'0.3631", "  -0.37%  ")
  10 = @("'1.269", "  +0.81%  ")
  11 = @("'1.001", "  +0.04%  ")
  12 = @("'0.08143", "  +0.13%  ")
  13 = @("'22.86", "  +1.50%  ")
  14 = @("'6.587", "  -0.03%  ")
  15 = @("'7.397", "  +0.65%  ")
  16 = @("'0.00001251", "  +0.47%  ")
  17 = @("1.606.28", "  -0.10%  ")
  18 = @("'94.01", "  +2.20%  ")
  19 = @("'0.06923", "  +1.34%  ")
  20 = @("'18.14", "  -0.44%  ")
  21 = @("'6.538", "  -0.04%  ")
  22 = @("'1.002", "  +0.04%  ")
  23 = @("'12.92", "  -1.13%  ")
  24 = @("23.235.02", "  +0.88%  ")
  25 = @("'2.448", "  +3.75%  ")
  26 = @("'3.076", "  +8.92%  ")
  27 = @("'21.17", "  +0.57%  ")
  28 = @("'150.08", "  -0.14%  ")
  29 = @("'5.275", "  +0.94%  ")
  30 = @("'135.33", "  +0.83%  ")
  31 = @("'2.386", "  +1.91%  ")
  32 = @("'6.762", "  -0.86%  ")
  33 = @("1.782.41", "  -0.17%  ")
  34 = @("'0.9646", "  +0.04%  ")
  35 = @("'0.07500", "  -0.95%  ")
  36 = @("'0.02761", "  +2.12%  ")
  37 = @("'10.35", "  +0.29%  ")
  38 = @("'0.2520", "  -0.28%  ")
  39 = @("'6.134", "  -1.91%  ")
  40 = @("'0.08803", "  -0.80%  ")
  41 = @("'1.419", "  +3.94%  ")
  42 = @("'0.7096", "  +0.95%  ")
  43 = @("'12.50", "  +0.98%  ")
  44 = @("'15.82", "  +4.14%  ")
  45 = @("'0.6544", "  -1.20%  ")
  46 = @("'2.328", "  +1.66%  ")
  47 = @("'4.007", "  +0.32%  ")
  48 = @("'132.96", "  +0.33%  ")
  49 = @("'0.07941", "  +0.55%  ")
  50 = @("'1.208", "  -1.36%  ")
  51 = @("'1.191", "  -2.97%  ")
}

foreach ($rowNum in $priceVolume.Keys) {
  $vals = $priceVolume[$rowNum]
  $dCell = $ws.Range("D$rowNum")
  $dCell.Value = $vals[0]
  $dCell.Style = "Normal"
  $eCell = $ws.Range("E$rowNum")
  $eCell.Value = $vals[1]
  $eCell.Style = "Normal"
}
